# Append the new "DSM Scheduled Flights vs actual" daily rows (970:991)
# that were collected after the previous update (through row 969).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the formatting (date / integer / percent styles) down from
#        the last existing data row (969) onto the new rows (970:991), the
#        same way Excel would if you dragged the fill handle down. ---
$ws.Range("A969:D969").Copy()
$ws.Range("A970:D991").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. New daily schedule data: date serial, scheduled (B), actual (C) ---
$rows = @(
  @(970, 44898, 48, 48),
  @(971, 44899, 63, 63),
  @(972, 44900, 66, 63),
  @(973, 44901, 75, 72),
  @(974, 44902, 74, 69),
  @(975, 44903, 83, 77),
  @(976, 44904, 66, 62),
  @(977, 44905, 50, 48),
  @(978, 44906, 60, 59),
  @(979, 44907, 68, 67),
  @(980, 44908, 63, 58),
  @(981, 44909, 65, 64),
  @(982, 44910, 72, 67),
  @(983, 44911, 78, 76),
  @(984, 44912, 52, 51),
  @(985, 44913, 55, 53),
  @(986, 44914, 56, 54),
  @(987, 44915, 76, 72),
  @(988, 44916, 73, 63),
  @(989, 44917, 62, 45),
  @(990, 44918, 57, 35),
  @(991, 44919, 53, 41)
)

foreach ($r in $rows) {
  $row = $r[0]
  $ws.Cells.Item($row, 1).Value = $r[1]
  $ws.Cells.Item($row, 2).Value = $r[2]
  $ws.Cells.Item($row, 3).Value = $r[3]
  $ws.Cells.Item($row, 4).Formula = "=C$row/B$row"
}

# --- 3. Scroll / selection state matches where the user ended up editing ---
$excel.ActiveWindow.ScrollRow = 953
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F984").Select() | Out-Null
